# COTB Study Test Suite
# Fix the study designation filter used by the "Programs/Studies/Cases..."
# summary query (column C, rows 2-5) from 'COT007B' to 'COTC007B'.
# (Column B formulas on each row reference their own distinct query text and
# must stay untouched - e.g. row 5's File query also contains 'COT007B' but
# is a different shared string and is not part of this fix.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Replace("COT007B", "COTC007B")

# Reflect the last active cell/selection as left after the edit.
$ws.Range("D5").Select()
